$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.387.73'
$ws.Range('E2').Value = '  -4.28%  '
$ws.Range('D3').Value = '2.981.19'
$ws.Range('E3').Value = '  -3.87%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '532.98'
$ws.Range('E5').Value = '  -1.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.38'
$ws.Range('E6').Value = '  -4.75%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '2.978.58'
$ws.Range('E8').Value = '  -3.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.493'
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.147'
$ws.Range('E10').Value = '  -6.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.05'
$ws.Range('E11').Value = '  -6.37%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.443'
$ws.Range('E12').Value = '  -3.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000220'
$ws.Range('E13').Value = '  -3.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.34'
$ws.Range('E14').Value = '  -4.69%  '
$ws.Range('D15').Value = '3.452.01'
$ws.Range('E15').Value = '  -4.06%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '61.386.57'
$ws.Range('E16').Value = '  -4.11%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.109'
$ws.Range('E17').Value = '  -2.91%  '
$ws.Range('D18').Value = '2.972.27'
$ws.Range('E18').Value = '  -3.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.55'
$ws.Range('E19').Value = '  -2.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '461.00'
$ws.Range('E20').Value = '  -5.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.07'
$ws.Range('E21').Value = '  -3.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.664'
$ws.Range('E22').Value = '  -6.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.88'
$ws.Range('E23').Value = '  -4.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.70'
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.77'
$ws.Range('E25').Value = '  -4.55%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.65'
$ws.Range('E27').Value = '  -2.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.62'
$ws.Range('E28').Value = '  -8.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('B30').Value = 'Mantle'
$ws.Range('C30').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.15'
$ws.Range('E30').Value = '  -2.04%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.33'
$ws.Range('E31').Value = '  -3.97%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.86'
$ws.Range('E32').Value = '  -2.57%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '54.87'
$ws.Range('E33').Value = '  -3.79%  '
$ws.Range('B34').Value = 'Stacks'
$ws.Range('C34').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.25'
$ws.Range('E34').Value = '  -6.46%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.31'
$ws.Range('E35').Value = '  -4.64%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.83'
$ws.Range('E36').Value = '  -4.22%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '455.94'
$ws.Range('E37').Value = '  -9.81%  '
$ws.Range('D38').Value = '3.126.82'
$ws.Range('E38').Value = '  -4.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0778'
$ws.Range('E39').Value = '  -2.92%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0378'
$ws.Range('E40').Value = '  -6.07%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.118'
$ws.Range('E41').Value = '  -1.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.01'
$ws.Range('E42').Value = '  -2.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.38'
$ws.Range('E43').Value = '  -12.18%  '
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '25.58'
$ws.Range('E45').Value = '  +1.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.240'
$ws.Range('E46').Value = '  -7.28%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.98'
$ws.Range('E47').Value = '  -7.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.107'
$ws.Range('E48').Value = '  -2.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '116.19'
$ws.Range('E49').Value = '  -5.04%  '
$ws.Range('D50').Value = '0.0₃0486'
$ws.Range('E50').Value = '  -10.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.25'
$ws.Range('E51').Value = '  +5.69%  '
